# Vendor Creation And test data
# Applies the TestData.xlsx edit: new "AP Vendor: New" page reference,
# a new AR sales-order test row, and two rows of invoice test data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Update existing AC9 label: "NetChain2 - AP Vendors" -> "NetChain2 - AP Vendor: New"
#        (this becomes shared-string index 61)
$ws.Range("AC9").Value = "NetChain2 – AP Vendor: New"

# --- 2. New row 13: AR sales order test entry
#        (adds shared strings 62 "AR.NetchainTest.CreateSalesOrder", 63 "Accenture")
$ws.Range("A13").Value = "AR.NetchainTest.CreateSalesOrder"
$ws.Range("B13").Value = "Accenture"
$ws.Range("C13").Value = "pune"
$ws.Rows.Item(13).RowHeight = 35.25

# --- 3. Fill out row 5 (invoice test data row)
#        (adds shared strings 64 "workbooks", 65 "measure2")
$ws.Range("B5").Value = "tech m"
$ws.Range("C5").Value = "Net 30"
$ws.Range("D5").Value = "pune"
$ws.Range("E5").Value = "Advertising"
$ws.Range("F5").Value = "invoice desc 1"
$ws.Range("G5").Value = 500
$ws.Range("H5").Value = "workbooks"
$ws.Range("I5").Value = "Department1"
$ws.Range("J5").Value = "Advertising"
$ws.Range("K5").Value = "inv desc2"
$ws.Range("L5").Value = "measure2"
$ws.Range("M5").Value = 234.44444999999999
$ws.Range("N5").Value = 876.9757366

# --- 4. Fill out row 6 (invoice test data row)
#        (adds shared strings 66 "Bank Charges", 67 "inv desc3", 68 "measure3")
$ws.Range("B6").Value = "tech m"
$ws.Range("C6").Value = "Net 30"
$ws.Range("D6").Value = "pune"
$ws.Range("E6").Value = "Advertising"
$ws.Range("F6").Value = "invoice desc 1"
$ws.Range("G6").Value = 500
$ws.Range("H6").Value = "laptop"
$ws.Range("I6").Value = "Department1"
$ws.Range("J6").Value = "Bank Charges"
$ws.Range("K6").Value = "inv desc3"
$ws.Range("L6").Value = "measure3"
$ws.Range("M6").Value = 38.741129999999998
$ws.Range("N6").Value = 34.987736650000002

# --- 5. Widen column AC to fit the new label text
$ws.Columns("AC").ColumnWidth = 26.1666666666667

# --- 6. Update the view: scroll/select so K6 is the active cell
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 1
$ws.Range("K6").Select()
